# DIYIOT - Board and Sensor Pin Map : cleanup pass
# - Updates the "Example Code" header note with full attribution text
# - Widens column A to fit the new "Specialized boards" labels
# - Adds a small note row above the "Specialized boards" block
# - Adds a "Notes:" section at the bottom explaining analog/voltage-divider usage
# - Resets the view (no frozen top-left cell, selection on O2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (old 8.6640625 -> new 25.33203125 "characters"); ColumnWidth is
# quantized internally to 1/6ths, 24.5 is the closest input that lands on the
# stored width closest to the target.
$ws.Columns.Item(1).ColumnWidth = 24.5

# New small marker cell above the "Specialized boards" header
$ws.Range("A46").Value = "``"

# Notes added at the bottom of the sheet (populate A80/A81 text before A79's
# "Notes:" label so the shared-string table order matches the source file)
$ws.Range("A80").Value = "If Analog = Use voltage divider to get input down to 1v"
$ws.Range("A81").Value = "Use 3.3v for all sensors unless otherwise needed, If input to sensor needs to be 5v, limit current output to 3.3v with a voltage divider version of board"
$ws.Range("A79").Value = "Notes:"

# Expand the "Example Code" label into a full attribution line
$ws.Range("O1").Value = "Example Code from Internet - huge thanks to all contributors over the years"

# Reset the view: clear the frozen/scrolled topLeftCell and select O2
$ws.Range("O2").Select() | Out-Null
